$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor-reading data (columns ax, ay, az, gx, gy, gz) for rows 2..21 after
# the sliding-window shift: two fresh rows were recorded and the two oldest
# windows were dropped (old row 22 removed), which is why the sheet's used
# range shrinks from A1:H22 to A1:H21.
$data = @(
    @(-3.823432922363281, 5.642979621887207, 4.475735664367676, -0.06746154287180267, 0.03214145534087062, -0.9956535080144557),
    @(-4.128349304199219, 5.92755126953125, 5.054780960083008, 0.2888078393284657, 0.04184071745636832, -1.758943790613218),
    @(-4.636096000671387, 5.381560325622559, 5.073278903961182, 1.293112953041876, -0.7586058728223629, -2.201892187429029),
    @(-5.365009784698486, 4.192014694213867, 5.306471347808838, 1.076085335986568, -0.3218408019043786, -1.197818033223937),
    @(-6.076753616333008, 2.580186367034912, 6.088143825531006, 1.466690946456983, -0.5813377291657031, -0.9997129280899739),
    @(-6.769103527069092, 3.8761305809021, 3.666451930999756, 4.411398765652687, -4.482220122980524, -3.045610698849662),
    @(-5.217415809631348, 5.673132419586182, 3.341332912445068, 0.08928575072178435, -1.550635254660315, 0.06088052863297655),
    @(-4.648830413818359, 2.288045167922974, 2.572498321533203, -0.9927981232487877, -1.707755399304766, 1.41453302322432),
    @(62.12903594970703, 9.40150260925293, -14.63701248168945, 2.431215341007984, 1.645354886387669, 0.273621466270713),
    @(-6.846994400024414, 9.073053359985352, 3.893637657165527, 0.2498580986677279, 0.9151578448539578, 0.4712292718332862),
    @(4.824018478393555, 1.011744022369385, 6.095863819122314, 0.04947298288691901, 0.4319243583568316, 0.3634221997371941),
    @(15.44786834716797, 3.326952695846558, 7.409719944000244, -0.5032464119410757, -0.05431018042009644, 0.3147802685582363),
    @(-0.8648982048034668, 5.278839111328125, 5.748791694641113, -0.8943670974221336, -1.390350583680838, 2.061716448429016),
    @(1.2574462890625, 1.225612640380859, 7.33539867401123, 0.08820253887841945, -0.5515473787347926, 0.8748975694179519),
    @(1.297415733337402, 4.380810260772705, 6.710616588592529, 0.2351689076181083, -0.03523484912029523, 0.06333464960199484),
    @(-4.824479103088379, 3.505787849426269, 5.917765617370605, -0.005455169480207168, -0.04037748242533482, 0.09177540264330601),
    @(2.601564407348633, 5.296027660369873, 4.190939426422119, -0.03705323929356984, -0.02731135879578294, 0.08058094122815343),
    @(0.9737215042114258, 3.227108478546143, 7.653890609741211, -0.03103692665003059, 0.001527163070128304, 0.01256180448500908),
    @(0.5169296264648438, 3.862967014312744, 6.639832496643066, -0.01491292483758091, -0.00055048926625147, -0.0008630249216113467),
    @(-0.4491424560546875, 3.653162002563477, 6.346769332885742, 0.02735042602343611, -0.0430127267920694, -0.02982940334220257)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Range("C$r").Value = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = $row[5]
}

# Remove the now-obsolete last row (was row 22); this also shrinks the
# worksheet's dimension/used range to A1:H21 automatically.
$ws.Rows.Item(22).Delete()
